# --------------------------------------------------------------------------
# rest_of_world_consolidated.xlsx -- regional LCA data cleanup
#   * lower-case / rename the header row to the new standard schema
#   * the data columns were re-derived upstream: the "Total MJ" column is
#     dropped and every remaining metric shifts one column to the left
#     (old Carbon -> new "carbon", old CED -> new "ced"); a brand-new
#     "climate change (kg CO2 eq)" column is appended in its place
#   * document each column's data type with a cell comment, as the
#     upstream data dictionary does
# --------------------------------------------------------------------------

$excel = New-Object -ComObject Excel.Application
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row: lower-case + rename to the new column schema -------------
$ws.Range("A1").Value = "industry"
$ws.Range("B1").Value = "unit"
$ws.Range("C1").Value = "process"
$ws.Range("D1").Value = "carbon (kg CO2 eq)"
$ws.Range("E1").Value = "ced (MJ)"
$ws.Range("F1").Value = "climate change (kg CO2 eq)"
$ws.Range("G1").Value = "region"

# --- data rows: shift D<-old E (carbon) and E<-old F (ced) for every row --
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 193) { $lastRow = 193 }
for ($r = 2; $r -le 193; $r++) {
    $oldCarbon = $ws.Cells.Item($r, 5).Value2
    $oldCed    = $ws.Cells.Item($r, 6).Value2
    $ws.Cells.Item($r, 4).Value = $oldCarbon
    $ws.Cells.Item($r, 5).Value = $oldCed
}

# --- new "climate change (kg CO2 eq)" column (column F) -------------------
    $ws.Cells.Item(2, 6).Value = 0.000000068486147
    $ws.Cells.Item(3, 6).Value = 0.0000055832315
    $ws.Cells.Item(4, 6).Value = 0.0000034502044
    $ws.Cells.Item(5, 6).Value = 0.0000015800225
    $ws.Cells.Item(6, 6).Value = 0.0000012655352
    $ws.Cells.Item(7, 6).Value = 0.0000052366095
    $ws.Cells.Item(8, 6).Value = 0.00000007332116299999999
    $ws.Cells.Item(9, 6).Value = 0.0000059587123
    $ws.Cells.Item(10, 6).Value = 0.0000058501022
    $ws.Cells.Item(11, 6).Value = 0.0000015963668
    $ws.Cells.Item(12, 6).Value = 0.0000065577001
    $ws.Cells.Item(13, 6).Value = 0.0000036206124
    $ws.Cells.Item(14, 6).Value = 0.0000024957723
    $ws.Cells.Item(15, 6).Value = 0.0000058873484
    $ws.Cells.Item(16, 6).Value = 0.0000052373063
    $ws.Cells.Item(17, 6).Value = 0.0000045146226
    $ws.Cells.Item(18, 6).Value = 0.0000053959431
    $ws.Cells.Item(19, 6).Value = 0.0000069372678
    $ws.Cells.Item(20, 6).Value = 0.0000051743527
    $ws.Cells.Item(21, 6).Value = 0.0000060856982
    $ws.Cells.Item(22, 6).Value = 0.0000064070601
    $ws.Cells.Item(23, 6).Value = 0.0000046855543
    $ws.Cells.Item(24, 6).Value = 0.0000013678482
    $ws.Cells.Item(25, 6).Value = 0.0000061600915
    $ws.Cells.Item(26, 6).Value = 0.0000061474213
    $ws.Cells.Item(27, 6).Value = 0.000000068707607
    $ws.Cells.Item(28, 6).Value = 0.000003932421
    $ws.Cells.Item(29, 6).Value = 0.0000048770989
    $ws.Cells.Item(30, 6).Value = 0.000008344941500000001
    $ws.Cells.Item(31, 6).Value = 0.00000069701298
    $ws.Cells.Item(32, 6).Value = 0.000005896152
    $ws.Cells.Item(33, 6).Value = 0.0000049352876
    $ws.Cells.Item(34, 6).Value = 0.0000024215397
    $ws.Cells.Item(35, 6).Value = 0.0000047499481
    $ws.Cells.Item(36, 6).Value = 0.0000022284384
    $ws.Cells.Item(37, 6).Value = 0.0000012787778
    $ws.Cells.Item(38, 6).Value = 0.0000059829496
    $ws.Cells.Item(39, 6).Value = 0.0000067420615
    $ws.Cells.Item(40, 6).Value = 0.000000067926242
    $ws.Cells.Item(41, 6).Value = 0.000000072200756
    $ws.Cells.Item(42, 6).Value = 0.0000067253576
    $ws.Cells.Item(43, 6).Value = 0.0000025111834
    $ws.Cells.Item(44, 6).Value = 0.0000048758718
    $ws.Cells.Item(45, 6).Value = 0.0000041732111
    $ws.Cells.Item(46, 6).Value = 0.0000022193113
    $ws.Cells.Item(47, 6).Value = 0.0000065474455
    $ws.Cells.Item(48, 6).Value = 0.0000070500425
    $ws.Cells.Item(49, 6).Value = 0.0000031412137
    $ws.Cells.Item(50, 6).Value = 0.00000033279135
    $ws.Cells.Item(51, 6).Value = 0.0000040245644
    $ws.Cells.Item(52, 6).Value = 0.0000064987044
    $ws.Cells.Item(53, 6).Value = 0.000000070852969
    $ws.Cells.Item(54, 6).Value = 0.0000065474455
    $ws.Cells.Item(55, 6).Value = 0.0000053205195
    $ws.Cells.Item(56, 6).Value = 0.0000058825079
    $ws.Cells.Item(57, 6).Value = 0.000000063412836
    $ws.Cells.Item(58, 6).Value = 0.000001314157
    $ws.Cells.Item(59, 6).Value = 0.0000049770074
    $ws.Cells.Item(60, 6).Value = 0.0000022178183
    $ws.Cells.Item(61, 6).Value = 0.0000041081956
    $ws.Cells.Item(62, 6).Value = 0.0000065455432
    $ws.Cells.Item(63, 6).Value = 0.00000060083403
    $ws.Cells.Item(64, 6).Value = 0.000000073553441
    $ws.Cells.Item(65, 6).Value = 0.0000035065575
    $ws.Cells.Item(66, 6).Value = 0.0000039639588
    $ws.Cells.Item(67, 6).Value = 0.0000025328147
    $ws.Cells.Item(68, 6).Value = 0.0000067982301
    $ws.Cells.Item(69, 6).Value = 0.000004177133
    $ws.Cells.Item(70, 6).Value = 0.0000041470222
    $ws.Cells.Item(71, 6).Value = 0.0000071303087
    $ws.Cells.Item(72, 6).Value = 0.0000013075301
    $ws.Cells.Item(73, 6).Value = 0.0000042950699
    $ws.Cells.Item(74, 6).Value = 0.0000014957589
    $ws.Cells.Item(75, 6).Value = 0.000006663979
    $ws.Cells.Item(76, 6).Value = 0.0000058752562
    $ws.Cells.Item(77, 6).Value = 0.0000026168259
    $ws.Cells.Item(78, 6).Value = 0.0000024378048
    $ws.Cells.Item(79, 6).Value = 0.0000071303087
    $ws.Cells.Item(80, 6).Value = 0.0000067982301
    $ws.Cells.Item(81, 6).Value = 0.000012720091
    $ws.Cells.Item(82, 6).Value = 0.0000033911145
    $ws.Cells.Item(83, 6).Value = 0.0000067864699
    $ws.Cells.Item(84, 6).Value = 0.000000040261261
    $ws.Cells.Item(85, 6).Value = 0.0000075856726
    $ws.Cells.Item(86, 6).Value = 0.0000064817347
    $ws.Cells.Item(87, 6).Value = 0.0000054168382
    $ws.Cells.Item(88, 6).Value = 0.000010356808
    $ws.Cells.Item(89, 6).Value = 0.0000051996596
    $ws.Cells.Item(90, 6).Value = 0.0000062203007
    $ws.Cells.Item(91, 6).Value = 0.0000040456441
    $ws.Cells.Item(92, 6).Value = 0.000004361072
    $ws.Cells.Item(93, 6).Value = 0.0000064788345
    $ws.Cells.Item(94, 6).Value = 0.00000045703081
    $ws.Cells.Item(95, 6).Value = 0.0000061887133
    $ws.Cells.Item(96, 6).Value = 0.0000072355872
    $ws.Cells.Item(97, 6).Value = 0.0000060533569
    $ws.Cells.Item(98, 6).Value = 0.0000014079291
    $ws.Cells.Item(99, 6).Value = 0.0000025512353
    $ws.Cells.Item(100, 6).Value = 0.000005926892
    $ws.Cells.Item(101, 6).Value = 0.00000006087501
    $ws.Cells.Item(102, 6).Value = 0.0000023922548
    $ws.Cells.Item(103, 6).Value = 0.000017246426
    $ws.Cells.Item(104, 6).Value = 0.0000031904116
    $ws.Cells.Item(105, 6).Value = 0.0000042586859
    $ws.Cells.Item(106, 6).Value = 0.00000044931357
    $ws.Cells.Item(107, 6).Value = 0.0000054908053
    $ws.Cells.Item(108, 6).Value = 0.0000067109333
    $ws.Cells.Item(109, 6).Value = 0.0000043235723
    $ws.Cells.Item(110, 6).Value = 0.0000050158361
    $ws.Cells.Item(111, 6).Value = 0.000005695726
    $ws.Cells.Item(112, 6).Value = 0.0000048263747
    $ws.Cells.Item(113, 6).Value = 0.0000057681451
    $ws.Cells.Item(114, 6).Value = 0.0000077918565
    $ws.Cells.Item(115, 6).Value = 0.0000036660355
    $ws.Cells.Item(116, 6).Value = 0.0000069372678
    $ws.Cells.Item(117, 6).Value = 0.0000066438733
    $ws.Cells.Item(118, 6).Value = 0.0000010903576
    $ws.Cells.Item(119, 6).Value = 0.0000030128223
    $ws.Cells.Item(120, 6).Value = 0.00000054019568
    $ws.Cells.Item(121, 6).Value = 0.0000061887133
    $ws.Cells.Item(122, 6).Value = 0.000000082965194
    $ws.Cells.Item(123, 6).Value = 0.0000058934433
    $ws.Cells.Item(124, 6).Value = 0.0000008363372
    $ws.Cells.Item(125, 6).Value = 0.0000023463963
    $ws.Cells.Item(126, 6).Value = 0.000010068596
    $ws.Cells.Item(127, 6).Value = 0.0000047228949
    $ws.Cells.Item(128, 6).Value = 0.0000037560874
    $ws.Cells.Item(129, 6).Value = 0.0000055707805
    $ws.Cells.Item(130, 6).Value = 0.000005579731
    $ws.Cells.Item(131, 6).Value = 0.0000041765412
    $ws.Cells.Item(132, 6).Value = 0.0000052618098
    $ws.Cells.Item(133, 6).Value = 0.0000014843537
    $ws.Cells.Item(134, 6).Value = 0.0000044724798
    $ws.Cells.Item(135, 6).Value = 0.00000005969247199999999
    $ws.Cells.Item(136, 6).Value = 0.0000023530246
    $ws.Cells.Item(137, 6).Value = 0.0000063738959
    $ws.Cells.Item(138, 6).Value = 0.0000067990379
    $ws.Cells.Item(139, 6).Value = 0.0000052776428
    $ws.Cells.Item(140, 6).Value = 0.0000040413017
    $ws.Cells.Item(141, 6).Value = 0.0000028943148
    $ws.Cells.Item(142, 6).Value = 0.0000066324593
    $ws.Cells.Item(143, 6).Value = 0.000006625385
    $ws.Cells.Item(144, 6).Value = 0.0000069372678
    $ws.Cells.Item(145, 6).Value = 0.0000062464408
    $ws.Cells.Item(146, 6).Value = 0.0000053205195
    $ws.Cells.Item(147, 6).Value = 0.0000042566592
    $ws.Cells.Item(148, 6).Value = 0.0000066258528
    $ws.Cells.Item(149, 6).Value = 0.0000056022468
    $ws.Cells.Item(150, 6).Value = 0.0000051190404
    $ws.Cells.Item(151, 6).Value = 0.0000056180955
    $ws.Cells.Item(152, 6).Value = 0.0000056109575
    $ws.Cells.Item(153, 6).Value = 0.00000042304768
    $ws.Cells.Item(154, 6).Value = 0.000004889469
    $ws.Cells.Item(155, 6).Value = 0.0000061887133
    $ws.Cells.Item(156, 6).Value = 0.0000058679519
    $ws.Cells.Item(157, 6).Value = 0.0000070684604
    $ws.Cells.Item(158, 6).Value = 0.0000040651087
    $ws.Cells.Item(159, 6).Value = 0.0000060127833
    $ws.Cells.Item(160, 6).Value = 0.0000049478539
    $ws.Cells.Item(161, 6).Value = 0.0000026061184
    $ws.Cells.Item(162, 6).Value = 0.0000033339179
    $ws.Cells.Item(163, 6).Value = 0.0000060773
    $ws.Cells.Item(164, 6).Value = 0.0000056212419
    $ws.Cells.Item(165, 6).Value = 0.00000097470265
    $ws.Cells.Item(166, 6).Value = 0.0000040942142
    $ws.Cells.Item(167, 6).Value = 0.0000049126982
    $ws.Cells.Item(168, 6).Value = 0.000013171955
    $ws.Cells.Item(169, 6).Value = 0.0000054268384
    $ws.Cells.Item(170, 6).Value = 0.0000050877959
    $ws.Cells.Item(171, 6).Value = 0.0000056540537
    $ws.Cells.Item(172, 6).Value = 0.00000462768
    $ws.Cells.Item(173, 6).Value = 0.0000056788994
    $ws.Cells.Item(174, 6).Value = 0.0000069372678
    $ws.Cells.Item(175, 6).Value = 0.00000013156036
    $ws.Cells.Item(176, 6).Value = 0.0000023290998
    $ws.Cells.Item(177, 6).Value = 0.0000038663465
    $ws.Cells.Item(178, 6).Value = 0.0000039507116
    $ws.Cells.Item(179, 6).Value = 0.0000067363202
    $ws.Cells.Item(180, 6).Value = 0.0000007101828
    $ws.Cells.Item(181, 6).Value = 0.0000051955685
    $ws.Cells.Item(182, 6).Value = 0.0000053179992
    $ws.Cells.Item(183, 6).Value = 0.0000018931736
    $ws.Cells.Item(184, 6).Value = 0.000004740675
    $ws.Cells.Item(185, 6).Value = 0.00000006994522599999999
    $ws.Cells.Item(186, 6).Value = 0.0000043236935
    $ws.Cells.Item(187, 6).Value = 0.0000060884934
    $ws.Cells.Item(188, 6).Value = 0.0000010389055
    $ws.Cells.Item(189, 6).Value = 0.0000029856299
    $ws.Cells.Item(190, 6).Value = 0.0000052831038
    $ws.Cells.Item(191, 6).Value = 0.0000069372678
    $ws.Cells.Item(192, 6).Value = 0.0000072123296
    $ws.Cells.Item(193, 6).Value = 0.0000043029719

# --- header comments describing each column's data type -------------------
$ws.Range("A1").AddComment("Data type: Categorical (text)")
$ws.Range("B1").AddComment("Data type: Various (e.g. kg, kWh)")
$ws.Range("C1").AddComment("Data type: Categorical (text)")
$ws.Range("D1").AddComment("Data type: Carbon footprint")
$ws.Range("E1").AddComment("Data type: Cumulative energy demand")
$ws.Range("F1").AddComment("Data type: Climate change impact")
$ws.Range("G1").AddComment("Data type: Categorical (text)")
